$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92
$ws.Cells.Item(92, 2).Value = 6754065
$ws.Cells.Item(92, 6).Value = 'Necaxa'
$ws.Cells.Item(92, 7).Value = 'Cruz Azul'
$ws.Cells.Item(92, 8).Value = 1
$ws.Cells.Item(92, 9).Value = 3
$ws.Cells.Item(92, 10).Value = 'A'
$ws.Cells.Item(92, 11).Value = 2.375
$ws.Cells.Item(92, 12).Value = 3.3
$ws.Cells.Item(92, 13).Value = 2.8
$ws.Cells.Item(92, 14).Value = 3.5
$ws.Cells.Item(92, 16).Value = 2.1
$ws.Cells.Item(92, 17).Value = 0.25
$ws.Cells.Item(92, 18).Value = 2
$ws.Cells.Item(92, 19).Value = 1.85
$ws.Cells.Item(92, 21).Value = 1.9
$ws.Cells.Item(92, 22).Value = 1.95
$ws.Cells.Item(92, 23).Value = -1
$ws.Cells.Item(92, 25).Value = 1.1
$ws.Cells.Item(92, 26).Value = -1
$ws.Cells.Item(92, 27).Value = 0.8500000000000001
$ws.Cells.Item(92, 28).Value = 0.8999999999999999

# Row 93
$ws.Cells.Item(93, 2).Value = 6754066
$ws.Cells.Item(93, 6).Value = 'Unam Pumas'
$ws.Cells.Item(93, 7).Value = 'Queretaro'
$ws.Cells.Item(93, 8).Value = 4
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 'H'
$ws.Cells.Item(93, 11).Value = 1.727
$ws.Cells.Item(93, 12).Value = 3.5
$ws.Cells.Item(93, 13).Value = 4.5
$ws.Cells.Item(93, 14).Value = 1.8
$ws.Cells.Item(93, 16).Value = 4.5
$ws.Cells.Item(93, 17).Value = -0.75
$ws.Cells.Item(93, 18).Value = 2.025
$ws.Cells.Item(93, 19).Value = 1.825
$ws.Cells.Item(93, 21).Value = 1.825
$ws.Cells.Item(93, 22).Value = 2.025
$ws.Cells.Item(93, 23).Value = 0.8
$ws.Cells.Item(93, 25).Value = -1
$ws.Cells.Item(93, 26).Value = 1.025
$ws.Cells.Item(93, 27).Value = -1
$ws.Cells.Item(93, 28).Value = 0.825

# Row 94
$ws.Cells.Item(94, 2).Value = 6754067
$ws.Cells.Item(94, 6).Value = 'Tigres UANL'
$ws.Cells.Item(94, 7).Value = 'Toluca'
$ws.Cells.Item(94, 9).Value = 2
$ws.Cells.Item(94, 10).Value = 'D'
$ws.Cells.Item(94, 11).Value = 1.8
$ws.Cells.Item(94, 12).Value = 3.3
$ws.Cells.Item(94, 13).Value = 4.333
$ws.Cells.Item(94, 14).Value = 1.533
$ws.Cells.Item(94, 15).Value = 4.2
$ws.Cells.Item(94, 16).Value = 6
$ws.Cells.Item(94, 17).Value = -1
$ws.Cells.Item(94, 18).Value = 1.925
$ws.Cells.Item(94, 19).Value = 1.925
$ws.Cells.Item(94, 21).Value = 2
$ws.Cells.Item(94, 22).Value = 1.85
$ws.Cells.Item(94, 23).Value = -1
$ws.Cells.Item(94, 24).Value = 3.2
$ws.Cells.Item(94, 26).Value = -1
$ws.Cells.Item(94, 27).Value = 0.925
$ws.Cells.Item(94, 28).Value = 1
$ws.Cells.Item(94, 29).Value = -1

# Row 95
$ws.Cells.Item(95, 2).Value = 7260442
$ws.Cells.Item(95, 6).Value = 'Santos Laguna'
$ws.Cells.Item(95, 7).Value = 'Tijuana'
$ws.Cells.Item(95, 9).Value = 1
$ws.Cells.Item(95, 10).Value = 'H'
$ws.Cells.Item(95, 11).Value = 1.75
$ws.Cells.Item(95, 12).Value = 3.6
$ws.Cells.Item(95, 13).Value = 4.2
$ws.Cells.Item(95, 14).Value = 1.65
$ws.Cells.Item(95, 15).Value = 4
$ws.Cells.Item(95, 16).Value = 4.75
$ws.Cells.Item(95, 17).Value = -0.75
$ws.Cells.Item(95, 18).Value = 1.8
$ws.Cells.Item(95, 19).Value = 2.05
$ws.Cells.Item(95, 21).Value = 1.85
$ws.Cells.Item(95, 22).Value = 2
$ws.Cells.Item(95, 23).Value = 0.6499999999999999
$ws.Cells.Item(95, 24).Value = -1
$ws.Cells.Item(95, 26).Value = 0.4
$ws.Cells.Item(95, 27).Value = -0.5
$ws.Cells.Item(95, 28).Value = 0
$ws.Cells.Item(95, 29).Value = -0

# Row 188
$ws.Cells.Item(188, 2).Value = 7612817
$ws.Cells.Item(188, 6).Value = 'Atletico San Luis'
$ws.Cells.Item(188, 7).Value = 'Tigres UANL'
$ws.Cells.Item(188, 9).Value = 2
$ws.Cells.Item(188, 10).Value = 'A'
$ws.Cells.Item(188, 11).Value = 3.4
$ws.Cells.Item(188, 12).Value = 3.6
$ws.Cells.Item(188, 13).Value = 2.05
$ws.Cells.Item(188, 14).Value = 3
$ws.Cells.Item(188, 15).Value = 3.3
$ws.Cells.Item(188, 16).Value = 2.4
$ws.Cells.Item(188, 17).Value = 0.25
$ws.Cells.Item(188, 18).Value = 1.775
$ws.Cells.Item(188, 19).Value = 2.1
$ws.Cells.Item(188, 20).Value = 2.25
$ws.Cells.Item(188, 21).Value = 1.8
$ws.Cells.Item(188, 22).Value = 2.05
$ws.Cells.Item(188, 24).Value = -1
$ws.Cells.Item(188, 25).Value = 1.4
$ws.Cells.Item(188, 27).Value = 1.1
$ws.Cells.Item(188, 28).Value = 0.8
$ws.Cells.Item(188, 29).Value = -1

# Row 189
$ws.Cells.Item(189, 2).Value = 7612818
$ws.Cells.Item(189, 6).Value = 'Monterrey'
$ws.Cells.Item(189, 7).Value = 'Queretaro'
$ws.Cells.Item(189, 9).Value = 1
$ws.Cells.Item(189, 10).Value = 'D'
$ws.Cells.Item(189, 11).Value = 1.444
$ws.Cells.Item(189, 12).Value = 4.2
$ws.Cells.Item(189, 13).Value = 7.5
$ws.Cells.Item(189, 14).Value = 1.4
$ws.Cells.Item(189, 15).Value = 4.333
$ws.Cells.Item(189, 16).Value = 8.5
$ws.Cells.Item(189, 17).Value = -1.25
$ws.Cells.Item(189, 18).Value = 1.925
$ws.Cells.Item(189, 19).Value = 1.925
$ws.Cells.Item(189, 20).Value = 2.75
$ws.Cells.Item(189, 21).Value = 2.025
$ws.Cells.Item(189, 22).Value = 1.825
$ws.Cells.Item(189, 24).Value = 3.333
$ws.Cells.Item(189, 25).Value = -1
$ws.Cells.Item(189, 27).Value = 0.925
$ws.Cells.Item(189, 28).Value = -1
$ws.Cells.Item(189, 29).Value = 0.825

# Row 193
$ws.Cells.Item(193, 2).Value = 7612811
$ws.Cells.Item(193, 6).Value = 'Leon'
$ws.Cells.Item(193, 7).Value = 'Santos Laguna'
$ws.Cells.Item(193, 8).Value = 3
$ws.Cells.Item(193, 9).Value = 2
$ws.Cells.Item(193, 11).Value = 1.833
$ws.Cells.Item(193, 12).Value = 3.75
$ws.Cells.Item(193, 13).Value = 4
$ws.Cells.Item(193, 14).Value = 1.8
$ws.Cells.Item(193, 15).Value = 4
$ws.Cells.Item(193, 16).Value = 4.2
$ws.Cells.Item(193, 17).Value = -0.75
$ws.Cells.Item(193, 18).Value = 2
$ws.Cells.Item(193, 19).Value = 1.85
$ws.Cells.Item(193, 20).Value = 3
$ws.Cells.Item(193, 21).Value = 1.975
$ws.Cells.Item(193, 22).Value = 1.875
$ws.Cells.Item(193, 23).Value = 0.8
$ws.Cells.Item(193, 26).Value = 0.5
$ws.Cells.Item(193, 27).Value = -0.5
$ws.Cells.Item(193, 28).Value = 0.9750000000000001
$ws.Cells.Item(193, 29).Value = -1

# Row 194
$ws.Cells.Item(194, 2).Value = 7612810
$ws.Cells.Item(194, 6).Value = 'Cruz Azul'
$ws.Cells.Item(194, 7).Value = 'Mazatlan FC'
$ws.Cells.Item(194, 8).Value = 2
$ws.Cells.Item(194, 9).Value = 1
$ws.Cells.Item(194, 11).Value = 1.727
$ws.Cells.Item(194, 12).Value = 3.9
$ws.Cells.Item(194, 13).Value = 4.333
$ws.Cells.Item(194, 14).Value = 1.5
$ws.Cells.Item(194, 15).Value = 4.2
$ws.Cells.Item(194, 16).Value = 5.5
$ws.Cells.Item(194, 17).Value = -1
$ws.Cells.Item(194, 18).Value = 1.85
$ws.Cells.Item(194, 19).Value = 2
$ws.Cells.Item(194, 20).Value = 2.75
$ws.Cells.Item(194, 21).Value = 1.8
$ws.Cells.Item(194, 22).Value = 2.05
$ws.Cells.Item(194, 23).Value = 0.5
$ws.Cells.Item(194, 26).Value = 0
$ws.Cells.Item(194, 27).Value = -0
$ws.Cells.Item(194, 28).Value = 0.4
$ws.Cells.Item(194, 29).Value = -0.5

# Row 237
$ws.Cells.Item(237, 2).Value = 7612867
$ws.Cells.Item(237, 6).Value = 'Club America'
$ws.Cells.Item(237, 7).Value = 'Mazatlan FC'
$ws.Cells.Item(237, 9).Value = 2
$ws.Cells.Item(237, 10).Value = 'D'
$ws.Cells.Item(237, 11).Value = 1.363
$ws.Cells.Item(237, 12).Value = 5
$ws.Cells.Item(237, 13).Value = 7.5
$ws.Cells.Item(237, 14).Value = 1.222
$ws.Cells.Item(237, 15).Value = 6.5
$ws.Cells.Item(237, 16).Value = 12
$ws.Cells.Item(237, 17).Value = -1.75
$ws.Cells.Item(237, 18).Value = 1.825
$ws.Cells.Item(237, 19).Value = 2.025
$ws.Cells.Item(237, 20).Value = 3.25
$ws.Cells.Item(237, 21).Value = 1.975
$ws.Cells.Item(237, 22).Value = 1.875
$ws.Cells.Item(237, 24).Value = 5.5
$ws.Cells.Item(237, 25).Value = -1
$ws.Cells.Item(237, 27).Value = 1.025
$ws.Cells.Item(237, 28).Value = 0.9750000000000001

# Row 238
$ws.Cells.Item(238, 2).Value = 7612866
$ws.Cells.Item(238, 6).Value = 'Leon'
$ws.Cells.Item(238, 7).Value = 'Cruz Azul'
$ws.Cells.Item(238, 9).Value = 3
$ws.Cells.Item(238, 10).Value = 'A'
$ws.Cells.Item(238, 11).Value = 2.5
$ws.Cells.Item(238, 12).Value = 3.4
$ws.Cells.Item(238, 13).Value = 2.7
$ws.Cells.Item(238, 14).Value = 2.8
$ws.Cells.Item(238, 15).Value = 3.6
$ws.Cells.Item(238, 16).Value = 2.375
$ws.Cells.Item(238, 17).Value = 0.25
$ws.Cells.Item(238, 18).Value = 1.75
$ws.Cells.Item(238, 19).Value = 2.05
$ws.Cells.Item(238, 20).Value = 2.75
$ws.Cells.Item(238, 21).Value = 1.85
$ws.Cells.Item(238, 22).Value = 2
$ws.Cells.Item(238, 24).Value = -1
$ws.Cells.Item(238, 25).Value = 1.375
$ws.Cells.Item(238, 27).Value = 1.05
$ws.Cells.Item(238, 28).Value = 0.8500000000000001

# Row 272
$ws.Cells.Item(272, 2).Value = 7612894
$ws.Cells.Item(272, 6).Value = 'Tigres UANL'
$ws.Cells.Item(272, 7).Value = 'Mazatlan FC'
$ws.Cells.Item(272, 8).Value = 5
$ws.Cells.Item(272, 11).Value = 1.4
$ws.Cells.Item(272, 12).Value = 4.5
$ws.Cells.Item(272, 13).Value = 6.5
$ws.Cells.Item(272, 14).Value = 1.615
$ws.Cells.Item(272, 15).Value = 4
$ws.Cells.Item(272, 18).Value = 2.05
$ws.Cells.Item(272, 19).Value = 1.8
$ws.Cells.Item(272, 20).Value = 2.75
$ws.Cells.Item(272, 21).Value = 1.925
$ws.Cells.Item(272, 22).Value = 1.925
$ws.Cells.Item(272, 23).Value = 0.615
$ws.Cells.Item(272, 26).Value = 1.05
$ws.Cells.Item(272, 27).Value = -1
$ws.Cells.Item(272, 28).Value = 0.925
$ws.Cells.Item(272, 29).Value = -1

# Row 273
$ws.Cells.Item(273, 2).Value = 7612892
$ws.Cells.Item(273, 6).Value = 'Leon'
$ws.Cells.Item(273, 7).Value = 'Puebla'
$ws.Cells.Item(273, 8).Value = 2
$ws.Cells.Item(273, 11).Value = 1.571
$ws.Cells.Item(273, 12).Value = 4
$ws.Cells.Item(273, 13).Value = 4.75
$ws.Cells.Item(273, 14).Value = 1.5
$ws.Cells.Item(273, 15).Value = 4.75
$ws.Cells.Item(273, 18).Value = 1.8
$ws.Cells.Item(273, 19).Value = 2.05
$ws.Cells.Item(273, 20).Value = 3.25
$ws.Cells.Item(273, 21).Value = 2.05
$ws.Cells.Item(273, 22).Value = 1.8
$ws.Cells.Item(273, 23).Value = 0.5
$ws.Cells.Item(273, 26).Value = 0
$ws.Cells.Item(273, 27).Value = -0
$ws.Cells.Item(273, 28).Value = -0.5
$ws.Cells.Item(273, 29).Value = 0.4

# Row 293
$ws.Cells.Item(293, 2).Value = 7612913
$ws.Cells.Item(293, 6).Value = 'Santos Laguna'
$ws.Cells.Item(293, 7).Value = 'Club America'
$ws.Cells.Item(293, 8).Value = 1
$ws.Cells.Item(293, 10).Value = 'D'
$ws.Cells.Item(293, 11).Value = 4.333
$ws.Cells.Item(293, 12).Value = 3.75
$ws.Cells.Item(293, 13).Value = 1.727
$ws.Cells.Item(293, 14).Value = 5.75
$ws.Cells.Item(293, 15).Value = 4
$ws.Cells.Item(293, 16).Value = 1.571
$ws.Cells.Item(293, 17).Value = 1
$ws.Cells.Item(293, 18).Value = 1.8
$ws.Cells.Item(293, 19).Value = 2.05
$ws.Cells.Item(293, 23).Value = -1
$ws.Cells.Item(293, 24).Value = 3
$ws.Cells.Item(293, 26).Value = 0.8
$ws.Cells.Item(293, 28).Value = -1
$ws.Cells.Item(293, 29).Value = 0.875

# Row 294
$ws.Cells.Item(294, 2).Value = 7612912
$ws.Cells.Item(294, 6).Value = 'Cruz Azul'
$ws.Cells.Item(294, 7).Value = 'Monterrey'
$ws.Cells.Item(294, 8).Value = 2
$ws.Cells.Item(294, 10).Value = 'H'
$ws.Cells.Item(294, 11).Value = 2.4
$ws.Cells.Item(294, 12).Value = 3.2
$ws.Cells.Item(294, 13).Value = 2.875
$ws.Cells.Item(294, 14).Value = 1.95
$ws.Cells.Item(294, 15).Value = 3.5
$ws.Cells.Item(294, 16).Value = 3.8
$ws.Cells.Item(294, 17).Value = -0.5
$ws.Cells.Item(294, 18).Value = 1.975
$ws.Cells.Item(294, 19).Value = 1.875
$ws.Cells.Item(294, 23).Value = 0.95
$ws.Cells.Item(294, 24).Value = -1
$ws.Cells.Item(294, 26).Value = 0.9750000000000001
$ws.Cells.Item(294, 28).Value = 0.4875
$ws.Cells.Item(294, 29).Value = -0.5

# Row 297
$ws.Cells.Item(297, 14).Value = 2.1
$ws.Cells.Item(297, 15).Value = 3.4
$ws.Cells.Item(297, 16).Value = 3.6

# Row 298
$ws.Cells.Item(298, 14).Value = 2.3
$ws.Cells.Item(298, 16).Value = 3.2

# Row 299
$ws.Cells.Item(299, 14).Value = 4.333

# Row 300
$ws.Cells.Item(300, 15).Value = 3.4

# Row 301
$ws.Cells.Item(301, 14).Value = 1.75
$ws.Cells.Item(301, 15).Value = 3.6
$ws.Cells.Item(301, 21).Value = 1.975
$ws.Cells.Item(301, 22).Value = 1.875

# Row 302
$ws.Cells.Item(302, 15).Value = 3.4
$ws.Cells.Item(302, 16).Value = 3.6

# Row 303
$ws.Cells.Item(303, 15).Value = 3.75
$ws.Cells.Item(303, 16).Value = 3.8

# Row 304
$ws.Cells.Item(304, 16).Value = 4

# Row 305
$ws.Cells.Item(305, 14).Value = 2.3
$ws.Cells.Item(305, 15).Value = 3.1
$ws.Cells.Item(305, 16).Value = 3.4
$ws.Cells.Item(305, 17).Value = -0.25
$ws.Cells.Item(305, 18).Value = 2.025
$ws.Cells.Item(305, 19).Value = 1.825
